$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percent number format matching the existing "Percent" style (0.0%) used
# throughout columns U:AI for rows 2, 4 and 6.
$pctFormat = "0.0%"

# ---------------------------------------------------------------------
# Row 2 (Adult / Black) - Totals denominator corrected from 103 to 100,
# and two new computed columns (X2 "above %", Y2 "ground %") added.
# ---------------------------------------------------------------------
$ws.Range("X2").Formula = "=29/100"
$ws.Range("X2").NumberFormat = $pctFormat
$ws.Range("Y2").Formula = "=71/100"
$ws.Range("Y2").NumberFormat = $pctFormat

$ws.Range("Z2").Formula = "=25/100"
$ws.Range("AA2").Formula = "=6/100"
$ws.Range("AB2").Formula = "=25/100"
$ws.Range("AC2").Formula = "=23/100"
$ws.Range("AD2").Formula = "=42/100"
$ws.Range("AE2").Formula = "=3/100"
$ws.Range("AF2").Formula = "=5/100"
$ws.Range("AH2").Formula = "=6/100"
$ws.Range("AI2").Formula = "=17/100"

# ---------------------------------------------------------------------
# Row 3 (Juvenile / Black) - new empty, percent-formatted placeholder
# cells for the two new columns, matching the other empty cells in
# this row (U3:W3, Z3:AI3).
# ---------------------------------------------------------------------
$ws.Range("X3").NumberFormat = $pctFormat
$ws.Range("Y3").NumberFormat = $pctFormat

# ---------------------------------------------------------------------
# Row 4 (Adult / Cinammon) - Totals denominator corrected from 392 to
# 384, and the two new computed columns added.
# ---------------------------------------------------------------------
$ws.Range("X4").Formula = "=93/384"
$ws.Range("X4").NumberFormat = $pctFormat
$ws.Range("Y4").Formula = "=282/384"
$ws.Range("Y4").NumberFormat = $pctFormat

$ws.Range("Z4").Formula = "=102/384"
$ws.Range("AA4").Formula = "=30/384"
$ws.Range("AB4").Formula = "=81/384"
$ws.Range("AC4").Formula = "=109/384"
$ws.Range("AD4").Formula = "=198/384"
$ws.Range("AE4").Formula = "=10/384"
$ws.Range("AF4").Formula = "=5/384"
$ws.Range("AH4").Formula = "=26/384"
$ws.Range("AI4").Formula = "=74/384"

# ---------------------------------------------------------------------
# Row 5 (Juvenile / Cinammon) - new empty, percent-formatted placeholder
# cells for the two new columns.
# ---------------------------------------------------------------------
$ws.Range("X5").NumberFormat = $pctFormat
$ws.Range("Y5").NumberFormat = $pctFormat

# ---------------------------------------------------------------------
# Row 6 (Adult / Gray) - Totals denominator corrected from 2468 to
# 2376, and the two new computed columns added.
# ---------------------------------------------------------------------
$ws.Range("X6").Formula = "=656/2376"
$ws.Range("X6").NumberFormat = $pctFormat
$ws.Range("Y6").Formula = "=1686/2376"
$ws.Range("Y6").NumberFormat = $pctFormat

$ws.Range("Z6").Formula = "=574/2376"
$ws.Range("AA6").Formula = "=225/2376"
$ws.Range("AB6").Formula = "=521/2376"
$ws.Range("AC6").Formula = "=592/2376"
$ws.Range("AD6").Formula = "=1144/2376"
$ws.Range("AE6").Formula = "=77/2376"
$ws.Range("AF6").Formula = "=34/2376"
$ws.Range("AG6").Formula = "=1/2376"
$ws.Range("AH6").Formula = "=116/2376"
$ws.Range("AI6").Formula = "=330/2376"

# ---------------------------------------------------------------------
# Update the active selection to match the saved view state.
# ---------------------------------------------------------------------
$ws.Range("AC11").Select()
